# Updated cryptos list on Fri Mar 29 05:25:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text happens to parse as a plain number need a
# leading apostrophe so Excel stores them as text (quote-prefixed) instead of
# silently converting them to a Number -- matching the column, which is
# entirely free-form text (e.g. "70.334.76", "3.58E+84" style values are not
# meant to be numeric). We collect those addresses and reset their style back
# to Normal afterwards so no stray quote-prefix style lingers on the cell.
$textFixRows = @()

# --- Row 17 / 18 swap: Uniswap <-> WrappedEther ---
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.580.84"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.72"
$textFixRows += "D18"
$ws.Range("E18").Value = "  +3.22%  "

# --- Row 38 / 39 / 40 rotation: Maker -> Bittensor -> InjectiveProtocol -> Maker ---
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'524.47"
$textFixRows += "D38"
$ws.Range("E38").Value = "  -4.36%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'37.69"
$textFixRows += "D39"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.641.83"
$ws.Range("E40").Value = "  +8.15%  "

# --- Price (D) / Volume (E) updates for all other rows ---
$ws.Range("D2").Value = "70.334.76"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "3.566.98"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("D4").Value = "'1.00"
$textFixRows += "D4"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'611.02"
$textFixRows += "D5"
$ws.Range("E5").Value = "  +4.05%  "

$ws.Range("D6").Value = "'187.25"
$textFixRows += "D6"
$ws.Range("E6").Value = "  +2.02%  "

$ws.Range("D7").Value = "3.562.65"
$ws.Range("E7").Value = "  +1.66%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "'0.215"
$textFixRows += "D10"
$ws.Range("E10").Value = "  +8.91%  "

$ws.Range("D11").Value = "'0.648"
$textFixRows += "D11"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("D12").Value = "'54.10"
$textFixRows += "D12"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").Value = "'9.46"
$textFixRows += "D14"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").Value = "4.134.77"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").Value = "70.433.57"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D19").Value = "'18.99"
$textFixRows += "D19"
$ws.Range("E19").Value = "  -1.51%  "

$ws.Range("D20").Value = "'573.25"
$textFixRows += "D20"
$ws.Range("E20").Value = "  +7.60%  "

$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").Value = "'0.996"
$textFixRows += "D22"
$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("D23").Value = "'17.45"
$textFixRows += "D23"
$ws.Range("E23").Value = "  -2.51%  "

$ws.Range("D24").Value = "'4.78"
$textFixRows += "D24"
$ws.Range("E24").Value = "  +4.79%  "

$ws.Range("D25").Value = "'4.90"
$textFixRows += "D25"
$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("D26").Value = "'93.87"
$textFixRows += "D26"
$ws.Range("E26").Value = "  -1.37%  "

$ws.Range("E27").Value = "  -0.96%  "

$ws.Range("D28").Value = "'10.94"
$textFixRows += "D28"
$ws.Range("E28").Value = "  -1.46%  "

$ws.Range("D29").Value = "'9.40"
$textFixRows += "D29"
$ws.Range("E29").Value = "  +3.75%  "

$ws.Range("D30").Value = "'32.37"
$textFixRows += "D30"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("D31").Value = "'7.08"
$textFixRows += "D31"
$ws.Range("E31").Value = "  -2.27%  "

$ws.Range("D32").Value = "'12.24"
$textFixRows += "D32"
$ws.Range("E32").Value = "  -1.10%  "

$ws.Range("E33").Value = "  +2.61%  "

$ws.Range("D34").Value = "'64.29"
$textFixRows += "D34"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").Value = "'3.73"
$textFixRows += "D35"
$ws.Range("E35").Value = "  +20.46%  "

$ws.Range("E36").Value = "  +2.87%  "

$ws.Range("E37").Value = "  -0.79%  "

$ws.Range("D41").Value = "'1.00"
$textFixRows += "D41"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "0.0₃0784"
$ws.Range("E42").Value = "  +3.20%  "

$ws.Range("D43").Value = "'3.52"
$textFixRows += "D43"
$ws.Range("E43").Value = "  +4.22%  "

$ws.Range("D44").Value = "'0.139"
$textFixRows += "D44"
$ws.Range("E44").Value = "  +3.20%  "

$ws.Range("D45").Value = "'0.0458"
$textFixRows += "D45"
$ws.Range("E45").Value = "  +4.63%  "

$ws.Range("D46").Value = "'3.51"
$textFixRows += "D46"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("D47").Value = "'2.97"
$textFixRows += "D47"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("E48").Value = "  +3.10%  "

$ws.Range("D49").Value = "'9.20"
$textFixRows += "D49"
$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("D51").Value = "'136.79"
$textFixRows += "D51"
$ws.Range("E51").Value = "  +0.08%  "

# Strip the quote-prefix style picked up above so the cell style matches
# the original (un-styled) Price cells exactly.
foreach ($addr in $textFixRows) {
    $ws.Range($addr).Style = "Normal"
}
